$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.5231316725978647
$ws1.Range("C2").Value = 0.08333333333333333
$ws1.Range("D2").Value = 0.8571428571428571
$ws1.Range("E2").Value = 0.1518987341772152
$ws1.Range("F2").Value = 0.3
$ws1.Range("G2").Value = 0.631578947368421
$ws1.Range("H2").Value = 0.6742241840556447
$ws1.Range("I2").Value = 24
$ws1.Range("J2").Value = 264
$ws1.Range("K2").Value = 270
$ws1.Range("L2").Value = 4

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")

# Row 2 (label "0")
$ws2.Range("B2").Value = 0.9854014598540146
$ws2.Range("C2").Value = 0.5056179775280899
$ws2.Range("D2").Value = 0.6683168316831684

# Row 3 (label "1")
$ws2.Range("B3").Value = 0.08333333333333333
$ws2.Range("C3").Value = 0.8571428571428571
$ws2.Range("D3").Value = 0.1518987341772152

# Row 4 (accuracy)
$ws2.Range("B4").Value = 0.5231316725978647
$ws2.Range("C4").Value = 0.5231316725978647
$ws2.Range("D4").Value = 0.5231316725978647
$ws2.Range("E4").Value = 0.5231316725978647

# Row 5 (macro avg)
$ws2.Range("B5").Value = 0.5343673965936739
$ws2.Range("C5").Value = 0.6813804173354735
$ws2.Range("D5").Value = 0.4101077829301918

# Row 6 (weighted avg)
$ws2.Range("B6").Value = 0.9404585638707778
$ws2.Range("C6").Value = 0.5231316725978647
$ws2.Range("D6").Value = 0.6425878161490639

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

# Row 2 (Actual 0)
$ws3.Range("B2").Value = 270
$ws3.Range("C2").Value = 264

# Row 3 (Actual 1)
$ws3.Range("B3").Value = 4
$ws3.Range("C3").Value = 24
